$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9.. down to 10..
$ws.Rows("9:9").Insert()

# New translation note string (goes into the freshly inserted row, columns A and B)
$newString = "変数203（ARGP攻撃種類）`n1斬　2打撃　3水　4火　5雷`n6誘惑　7食べ物`n特殊206（個別攻撃種類）`n1リンゴ　2皿　3卵"
$ws.Range("A9").Value = $newString
$ws.Range("B9").Value = $newString

# Mirror column A into column B for every row, except row 10
# (row 10 = "ネズミが通れるサイズの穴が開いている・・・", which stays A-only)
for ($r = 1; $r -le 54; $r++) {
    if ($r -eq 10) { continue }
    $aText = $ws.Cells.Item($r, 1).Text
    $ws.Cells.Item($r, 2).Value = $aText
}

# Remove now-stray extra cells that are no longer part of the layout
$ws.Range("C5").ClearContents()
$ws.Range("D47").ClearContents()
